$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. "Legends:" header - reuses the same bold/underline big-title look that
#    the report header (A1 "iRipple, Inc.") already uses.
# ---------------------------------------------------------------------------
$ws.Range("E24:P24").Merge()
$ws.Range("E24").Value = "Legends:"
$ws.Range("E24:P24").Font.Name = "Arial"
$ws.Range("E24:P24").Font.Size = 15
$ws.Range("E24:P24").Font.Bold = $true
$ws.Range("E24:P24").Font.Underline = $true

# ---------------------------------------------------------------------------
# 2. Legend entry #1 - blue swatch + "has request/remark" explanation.
# ---------------------------------------------------------------------------
$ws.Range("E25:E26").Merge()
$ws.Range("F25:P26").Merge()
$ws.Range("E25").Interior.Color = 13411113
$ws.Range("F25").Value = "Employee has request(s)/remark(s) for that day.`n*May incur late and/or undertime depending on his or her time-in and time-out."
$ws.Range("F25:P26").Font.Name = "Arial"
$ws.Range("F25:P26").Font.Size = 11
$ws.Range("F25:P26").Font.Bold = $true
$ws.Range("F25:P26").Font.Underline = $true

# ---------------------------------------------------------------------------
# 3. Legend entry #2 - yellow swatch + "half-day" explanation.
# ---------------------------------------------------------------------------
$ws.Range("E27:E28").Merge()
$ws.Range("F27:P28").Merge()
$ws.Range("E27").Interior.Color = 6737151
$ws.Range("F27").Value = "Employee is considered half-day because of his time-in or time-out."
$ws.Range("F27:P28").Font.Name = "Arial"
$ws.Range("F27:P28").Font.Size = 11
$ws.Range("F27:P28").Font.Bold = $true
$ws.Range("F27:P28").Font.Underline = $true

# ---------------------------------------------------------------------------
# 4. Legend entry #3 - red swatch + "absent" explanation.
# ---------------------------------------------------------------------------
$ws.Range("E29:E30").Merge()
$ws.Range("F29:P30").Merge()
$ws.Range("E29").Interior.Color = 6184671
$ws.Range("F29").Value = "Employee has no time-in and therefore, considered as absent."
$ws.Range("F29:P30").Font.Name = "Arial"
$ws.Range("F29:P30").Font.Size = 11
$ws.Range("F29:P30").Font.Bold = $true
$ws.Range("F29:P30").Font.Underline = $true

Write-Output "done"
